# Update gh-pages output generated at 456a3b4
# Applies the scraped-data refresh to 北京-漫展信息.xlsx:
#   - bump several "want-to-go" counts (column F) across all four sheets
#   - insert two brand-new 演出 (show) rows (2024-12-24, 2024-12-30) ahead
#     of the existing 2025-01-01 / 2025-01-04 rows on the 演出 sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions) - column F (想去人数) bumps
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$expoF = @{
    4  = 5985
    6  = 65
    12 = 687
    13 = 1608
    15 = 1632
    16 = 565
    17 = 193
    18 = 649
    19 = 4646
    20 = 102
    21 = 47
    22 = 676
    23 = 3362
    24 = 825
    27 = 14
    28 = 2340
    33 = 1246
    38 = 1270
    39 = 1247
    40 = 84
}
foreach ($row in $expoF.Keys) {
    $wsExpo.Range("F$row").Value = $expoF[$row]
}

# ---------------------------------------------------------------------
# Sheet "演出" (shows) - column F bumps for existing rows that are NOT
# shifted by the insert below (rows 6 and 20 stay put; the insert only
# affects rows >= 25).
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F6").Value = 109
$wsShow.Range("F20").Value = 312

# Insert two new rows at 25-26; existing rows 25 (2025-01-01 Beyond gig)
# and 26 (2025-01-04 concert) shift down to 27 and 28 automatically.
$wsShow.Rows("25:26").Insert()

# New row 25: 2024-12-24 "永恒乐队" tribute show
$wsShow.Range("A25").Value = 24
$wsShow.Range("B25:E25").NumberFormat = "@"
$wsShow.Range("B25").Value = "2024-12-24"
$wsShow.Range("C25").Value = "北京·【限时早鸟5折】 致敬beyond 「永恒乐队」平安夜专场·爱与和平经典金曲演唱会"
$wsShow.Range("D25").Value = "北京西城区西直门外大街135号（北京展览馆后） 北京展览馆剧场"
$wsShow.Range("E25").Value = "2024.12.24 19:30-12.24 21:00"
$wsShow.Range("F25").Value = 0
$wsShow.Range("G25").Value = 50
$wsShow.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=93647"
$wsShow.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202410/Gww1N7bN1729155447572.jpeg"

# New row 26: 2024-12-30 Chen Jia / Teresa Teng tribute show
$wsShow.Range("A26").Value = 25
$wsShow.Range("B26:E26").NumberFormat = "@"
$wsShow.Range("B26").Value = "2024-12-30"
$wsShow.Range("C26").Value = "北京·早鸟95折甜蜜蜜—陈佳2025邓丽君经典金曲新年专场演唱会"
$wsShow.Range("D26").Value = "复内大街49号（民族宫饭店旁边） 京演民族文化宫大剧院"
$wsShow.Range("E26").Value = "2024.12.30 19:30-12.30 21:30"
$wsShow.Range("F26").Value = 0
$wsShow.Range("G26").Value = 456
$wsShow.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=93643"
$wsShow.Range("I26").Value = "//i2.hdslb.com/bfs/openplatform/202410/A3PcE0TX1729152716857.jpeg"

# The shifted-down rows (old 25 -> 27, old 26 -> 28) keep their own data but
# their running index in column A advances by 2 to stay in sequence.
$wsShow.Range("A27").Value = 26
$wsShow.Range("A28").Value = 27

# ---------------------------------------------------------------------
# Sheet "本地生活" (local life) - column F bumps
# ---------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 751
$wsLocal.Range("F4").Value = 206
$wsLocal.Range("F5").Value = 309

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types) - column F bumps
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$allF = @{
    6  = 751
    7  = 206
    8  = 5985
    10 = 109
    22 = 1608
    24 = 1632
    25 = 565
    26 = 193
    27 = 649
    28 = 4646
    29 = 676
    30 = 3362
    31 = 825
    34 = 2340
    42 = 502
    49 = 84
}
foreach ($row in $allF.Keys) {
    $wsAll.Range("F$row").Value = $allF[$row]
}
